# Apply changes described by the commit:
#  - growth sheet: append an extra closing paren to the "B" helper-formulas
#    (A#&"(" -> A#&"()") for every row that originally produced a trailing "("
#  - growth sheet: widen column B, change active cell / selection, and make
#    "growth" the active (selected) sheet/tab in the workbook
#  - icons sheet: no longer the selected tab (tabSelected removed)

$wb = $excel.ActiveWorkbook
$growth = $wb.Worksheets.Item("growth")
$icons = $wb.Worksheets.Item("icons")

# Rows on "growth" whose column B formula currently ends in the literal "("
# and needs an extra ")" appended, turning e.g. =A3&"(" into =A3&"()"
$rowsWithParen = @(3,4,5,6,9,10,13,15,17,18,19,21,22,24,25,26,27,28,29,30,31,32)

foreach ($r in $rowsWithParen) {
    $cell = $growth.Range("B$r")
    $cell.Formula = '=A' + $r + '&"()"'
}

# Widen column B on growth to match the new, longer values
$growth.Columns.Item(2).ColumnWidth = 20.65

# Make "growth" the active sheet/tab and set its selection to I21.
# (Leaving "icons" untouched keeps its own existing selection state, while
# it naturally loses tabSelected since it is no longer the active sheet.)
$growth.Activate()
$growth.Range("I21").Select()
